# "Generate Report for Handoff"
#
# The b.md row has moved from "Handed back: in sync with en-US" to
# "Ready for handoff" on the Overview sheet, and on each locale sheet
# (zh-cn / de-de) its handoff file / handoff datetime now point at the
# newly generated xlf package instead of the old one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md row. Columns B (zh-cn) and C (de-de)
# both flip from "Handed back: in sync with en-US" to "Ready for handoff".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md row.
#   B3 Status                -> "Ready for handoff"
#   C3 Latest Handoff File   -> new xlf file (display text + hyperlink)
#   D3 Latest Handoff Datetime -> new timestamp
# The other hyperlinks on the sheet are untouched, but since hyperlink
# objects loaded from the file can't be edited in place, every hyperlink
# is recreated (delete-all then re-add) preserving the original target
# addresses for the unaffected cells.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-08 05:07:38"

$zhcn.Range("A1").Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9ace12dc4906ecb0f27ea3eefc46c6a4c14a4d14/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1f076b853f362e63023677d22ab349a307d9fe4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/68ad440fbc0091f43c020c507de6527f8b4953d9/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/811fde2e2fd3837cc347b552c7dd2ec4a8c0b7f1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9ace12dc4906ecb0f27ea3eefc46c6a4c14a4d14/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1f076b853f362e63023677d22ab349a307d9fe4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/68ad440fbc0091f43c020c507de6527f8b4953d9/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/811fde2e2fd3837cc347b552c7dd2ec4a8c0b7f1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9ace12dc4906ecb0f27ea3eefc46c6a4c14a4d14/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md row, same shape of edit as zh-cn.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-08 05:07:48"

$dede.Range("A1").Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9ace12dc4906ecb0f27ea3eefc46c6a4c14a4d14/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b5e32baaf88df5626386f970b882740f89362266/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8301a441470987a56f8e7b88324de21b9ebb57f9/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f5a02fbd8d47ab2d8372d7c5864b0845bfa936e5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9ace12dc4906ecb0f27ea3eefc46c6a4c14a4d14/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b5e32baaf88df5626386f970b882740f89362266/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8301a441470987a56f8e7b88324de21b9ebb57f9/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f5a02fbd8d47ab2d8372d7c5864b0845bfa936e5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9ace12dc4906ecb0f27ea3eefc46c6a4c14a4d14/.localization-config", "", "", ".localization-config")

Write-Output "Report generated for handoff: b.md is ready for handoff (zh-cn, de-de)."
